$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (values refreshed by GitHub Actions run).
# Cells whose new value looks like a plain number (e.g. "580.81") need their
# NumberFormat forced to Text first, otherwise Excel auto-converts the literal
# into a floating point number (losing exact text such as trailing zeros).

$ws.Range("D2").Value = '62.889.42'
$ws.Range("E2").Value = '  +5.98%  '
$ws.Range("D3").Value = '2.445.15'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.81'
$ws.Range("E5").Value = '  +4.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.84'
$ws.Range("E6").Value = '  +6.56%  '
$ws.Range("E8").Value = '  +2.01%  '
$ws.Range("D9").Value = '2.443.58'
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("E10").Value = '  +6.18%  '
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("E12").Value = '  +3.71%  '
$ws.Range("E13").Value = '  +5.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.01'
$ws.Range("E14").Value = '  +6.59%  '
$ws.Range("E15").Value = '  +9.37%  '
$ws.Range("D16").Value = '2.890.20'
$ws.Range("E16").Value = '  +3.91%  '
$ws.Range("D17").Value = '62.679.36'
$ws.Range("E17").Value = '  +5.72%  '
$ws.Range("D18").Value = '2.437.93'
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.95'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  +5.18%  '
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.05'
$ws.Range("E23").Value = '  +14.42%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.76'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '615.92'
$ws.Range("E26").Value = '  +12.06%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.37'
$ws.Range("E27").Value = '  +4.33%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0980'
$ws.Range("E28").Value = '  +8.11%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.565.06'
$ws.Range("E29").Value = '  +3.81%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.13'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("E31").Value = '  +9.59%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.87'
$ws.Range("E32").Value = '  +6.01%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.137'
$ws.Range("E33").Value = '  +5.69%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.48'
$ws.Range("E34").Value = '  +5.86%  '
$ws.Range("D35").Value = '0.0₆0372'
$ws.Range("E35").Value = '  +30.70%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.75'
$ws.Range("E37").Value = '  +5.47%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.372'
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '152.61'
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.38'
$ws.Range("E40").Value = '  +7.94%  '
$ws.Range("B41").Value = 'EthereumClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.60'
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  +17.20%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.76'
$ws.Range("E43").Value = '  +7.78%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.37'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.98'
$ws.Range("E46").Value = '  +4.23%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.57'
$ws.Range("E47").Value = '  +2.61%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.15'
$ws.Range("E48").Value = '  +6.27%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.597'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0513'
$ws.Range("E50").Value = '  +3.50%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0916'
$ws.Range("E51").Value = '  +2.95%  '
